# Koodia siistitty ja kommentoitu, poistettu ääkköset taulukosta.
# Replace names that contain Finnish umlauts (ä/ö) with plain-ASCII
# equivalents (a/o) in the board-member list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old (with umlauts) -> new (umlauts stripped) names.
$replacements = @{
    "Arvi Syrjänen"  = "Arvi Syrjanen"
    "Ella Seppä"     = "Ella Seppa"
    "Jari Leppänen"  = "Jari Leppanen"
    "Paula Jyrkönen" = "Paula Jyrkonen"
}

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $val = $cell.Value()
    if ($replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}

# Restore the view so it scrolls to show row 16 onward, with A20 selected
# (matches the state the author left the sheet in before committing).
$ws.Activate()
$ws.Range("A20").Select()
$excel.ActiveWindow.ScrollRow = 16
